$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "March 5th" to "work"
$ws.Name = "work"

# Remove the stray date value in A3 (row 3 now has no date, only the other columns)
$ws.Range("A3").Clear()

# Add new row 14 for the new "Show nearest k nodes besides accident spot" task
$ws.Range("A14").Value = 43901
$ws.Range("B14").Value = "1"
$ws.Range("C14").Value = "New feature: Show nearest k nodes besides the accident spot."
$ws.Range("D14").Value = 0.60416666666666663
$ws.Range("E14").Value = 0.72916666666666663
$ws.Range("F14").Value = 3

# Update the selected cell to reflect the new last-used area
$ws.Range("C15").Select()
